$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J - copy formatting (style) from H1,
# which already carries the bold/centered/bordered header style, then set text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I and J values for data rows 2..71 (row -> I,J)
$values = @{
    2  = @(8,8);   3  = @(7,7);   4  = @(8,8);   5  = @(8,8);   6  = @(10,10);
    7  = @(8,8);   8  = @(9,9);   9  = @(8,8);   10 = @(8,8);   11 = @(8,8);
    12 = @(8,8);   13 = @(7,7);   14 = @(8,8);   15 = @(9,9);   16 = @(7,8);
    17 = @(7,7);   18 = @(9,9);   19 = @(8,8);   20 = @(8,8);   21 = @(8,8);
    22 = @(8,8);   23 = @(8,8);   24 = @(9,9);   25 = @(8,8);   26 = @(8,8);
    27 = @(8,8);   28 = @(8,8);   29 = @(8,8);   30 = @(8,8);   31 = @(8,8);
    32 = @(7,7);   33 = @(9,9);   34 = @(9,9);   35 = @(9,9);   36 = @(8,8);
    37 = @(9,9);   38 = @(10,10); 39 = @(8,8);   40 = @(9,9);   41 = @(10,10);
    42 = @(8,8);   43 = @(8,8);   44 = @(8,8);   45 = @(8,8);   46 = @(8,8);
    47 = @(7,7);   48 = @(8,8);   49 = @(8,8);   50 = @(8,8);   51 = @(8,8);
    52 = @(6,7);   53 = @(9,9);   54 = @(8,8);   55 = @(7,7);   56 = @(8,8);
    57 = @(8,8);   58 = @(8,9);   59 = @(8,8);   60 = @(7,8);   61 = @(7,7);
    62 = @(8,8);   63 = @(8,8);   64 = @(8,8);   65 = @(7,7);   66 = @(6,7);
    67 = @(4,4);   68 = @(8,8);   69 = @(5,5);   70 = @(4,4);   71 = @(4,4);
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
